$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain Text (they contain plain numeric-looking values;
# without this Excel auto-converts them to Number on assignment)
$textCells = @("D5","D6","D9","D10","D14","D15","D16","D19","D21","D22","D23","D24","D25","D27","D29","D30","D32","D33","D38","D39","D41","D43","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values per the latest scrape
$ws.Range('D2').Value = '27.101.59'
$ws.Range('E2').Value = '  -0.73%  '

$ws.Range('D3').Value = '1.629.65'
$ws.Range('E3').Value = '  -1.42%  '

$ws.Range('E4').Value = '  +0.34%  '

$ws.Range('D5').Value = '215.91'
$ws.Range('E5').Value = '  -1.41%  '

$ws.Range('D6').Value = '0.516'
$ws.Range('E6').Value = '  +1.22%  '

$ws.Range('E7').Value = '  +0.33%  '

$ws.Range('E8').Value = '  -1.46%  '

$ws.Range('D9').Value = '0.0623'
$ws.Range('E9').Value = '  -0.91%  '

$ws.Range('D10').Value = '20.03'
$ws.Range('E10').Value = '  -1.18%  '

$ws.Range('E11').Value = '  +0.05%  '

$ws.Range('D12').Value = '1.865.94'
$ws.Range('E12').Value = '  -0.89%  '

$ws.Range('D13').Value = '1.618.12'
$ws.Range('E13').Value = '  -1.49%  '

$ws.Range('D14').Value = '4.11'
$ws.Range('E14').Value = '  -0.94%  '

$ws.Range('D15').Value = '0.541'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('D16').Value = '65.85'
$ws.Range('E16').Value = '  -3.14%  '

$ws.Range('D17').Value = '27.089.95'
$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  -1.07%  '

$ws.Range('D19').Value = '214.08'
$ws.Range('E19').Value = '  -3.45%  '

$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').Value = '6.82'
$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('B22').Value = 'Toncoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D22').Value = '2.51'
$ws.Range('E22').Value = '  +0.70%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '4.38'
$ws.Range('E23').Value = '  -1.74%  '

$ws.Range('D24').Value = '9.08'
$ws.Range('E24').Value = '  -2.23%  '

$ws.Range('D25').Value = '147.16'
$ws.Range('E25').Value = '  -0.54%  '

$ws.Range('E26').Value = '  +0.46%  '

$ws.Range('D27').Value = '7.36'
$ws.Range('E27').Value = '  -0.70%  '

$ws.Range('E28').Value = '  -1.33%  '

$ws.Range('D29').Value = '15.57'
$ws.Range('E29').Value = '  -2.03%  '

$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  -0.65%  '

$ws.Range('E31').Value = '  -0.93%  '

$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  -0.70%  '

$ws.Range('D33').Value = '3.01'
$ws.Range('E33').Value = '  -0.98%  '

$ws.Range('D34').Value = '1.301.59'
$ws.Range('E34').Value = '  +2.52%  '

$ws.Range('E35').Value = '  -2.14%  '

$ws.Range('E36').Value = '  +0.48%  '

$ws.Range('E37').Value = '  -1.42%  '

$ws.Range('D38').Value = '0.542'
$ws.Range('E38').Value = '  -0.33%  '

$ws.Range('D39').Value = '0.845'
$ws.Range('E39').Value = '  -0.15%  '

$ws.Range('E40').Value = '  +0.31%  '

$ws.Range('D41').Value = '2.26'
$ws.Range('E41').Value = '  +3.31%  '

$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('D43').Value = '5.31'
$ws.Range('E43').Value = '  -1.56%  '

$ws.Range('D44').Value = '1.771.40'
$ws.Range('E44').Value = '  -1.18%  '

$ws.Range('E45').Value = '  -2.22%  '

$ws.Range('D46').Value = '90.37'
$ws.Range('E46').Value = '  -2.47%  '

$ws.Range('D47').Value = '1.59'
$ws.Range('E47').Value = '  -0.92%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '0.786'
$ws.Range('E48').Value = '  +17.48%  '

$ws.Range('D49').Value = '0.0512'
$ws.Range('E49').Value = '  -0.47%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0987'
$ws.Range('E50').Value = '  -6.89%  '

$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  -2.12%  '
